$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with revised totals (no row shift)
$ws.Range("C36").Value = "'793"
$ws.Range("E36").Value = "'3798794.93"

$ws.Range("C56").Value = "'1025"
$ws.Range("E56").Value = "'5819203.93"

$ws.Range("C57").Value = "'507"
$ws.Range("E57").Value = "'4652636.78"

$ws.Range("C58").Value = "'188"
$ws.Range("E58").Value = "'2053986.11"

# Update rows 76-79 (Martinique) with revised totals (no row shift, same categories)
$ws.Range("C76").Value = "'36"
$ws.Range("D76").Value = "'35"
$ws.Range("E76").Value = "'93074.00"

$ws.Range("C77").Value = "'230"
$ws.Range("D77").Value = "'226"
$ws.Range("E77").Value = "'701456.34"

$ws.Range("C78").Value = "'112"
$ws.Range("D78").Value = "'111"
$ws.Range("E78").Value = "'408140.60"

$ws.Range("C79").Value = "'36"
$ws.Range("D79").Value = "'36"
$ws.Range("E79").Value = "'155000.00"

# Insert a new row at position 80 for the new "Martinique / 10 a 19 salaries" category,
# shifting old row 80 (and everything below) down by one
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80
$ws.Range("A80").Value = "'Fonds de solidarité"
$ws.Range("B80").Value = "'VOLET2"
$ws.Range("C80").Value = "'3"
$ws.Range("D80").Value = "'3"
$ws.Range("E80").Value = "'22000.00"
$ws.Range("F80").Value = "'02"
$ws.Range("G80").Value = "'Martinique"
$ws.Range("H80").Value = "'11"
$ws.Range("I80").Value = "'10 à 19 salariés"

# Row 81 now holds the old row-80 data (shifted down); refresh its revised totals
$ws.Range("C81").Value = "'24"
$ws.Range("D81").Value = "'23"
$ws.Range("E81").Value = "'48150.00"
